$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Rename the month tabs to drop accented characters
# ------------------------------------------------------------------
$wb.Worksheets.Item("Mês - Janeiro").Name   = "Mes - Janeiro"
$wb.Worksheets.Item("Mês - Fevereiro").Name = "Mes - Fevereiro"
$wb.Worksheets.Item("Mês - Março").Name     = "Mes - Marco"
$wb.Worksheets.Item("Mês - Abril").Name     = "Mes - Abril"
$wb.Worksheets.Item("Mês - Maio").Name      = "Mes - Maio"
$wb.Worksheets.Item("Mês - Julho").Name     = "Mes - Julho"

# ------------------------------------------------------------------
# 2) "Geral" sheet - update partial scores (column B) for rows 2-21
# ------------------------------------------------------------------
$geral = $wb.Worksheets.Item("Geral")

$geralScores = @{
    2  = 51.5
    3  = 62.76
    4  = 38.46
    5  = 58.17
    6  = 83.5
    7  = 45.46
    8  = 63.76
    9  = 44.65
    10 = 53.06
    11 = 53.91
    12 = 39.66
    13 = 44.26
    14 = 54.36
    15 = 55.96
    17 = 16.4
    18 = 52.66
    19 = 67.16
    20 = 40.4
    21 = 59.86
}

foreach ($row in $geralScores.Keys) {
    $geral.Cells.Item($row, 2).Value = $geralScores[$row]
}

# Add new row 22 - "Lider_Rodada" totals row, copying the look of row 21
$geral.Range("A21:AM21").Copy()
$geral.Range("A22:AM22").PasteSpecial(-4122)
$geral.Range("A22").Value = "Lider_Rodada"
for ($c = 2; $c -le 39; $c++) {
    $geral.Cells.Item(22, $c).Value = 0
}

# ------------------------------------------------------------------
# 3) "Turno 2" sheet - add the matching new row 22
# ------------------------------------------------------------------
$turno2 = $wb.Worksheets.Item("Turno 2")

$turno2.Range("A21:T21").Copy()
$turno2.Range("A22:T22").PasteSpecial(-4122)
$turno2.Range("A22").Value = "Lider_Rodada"
for ($c = 2; $c -le 20; $c++) {
    $turno2.Cells.Item(22, $c).Value = 0
}

# ------------------------------------------------------------------
# 4) "Classif Turno 2" sheet - re-rank teams by the new partial scores
#    (descending order, matching the "Geral" sheet column B values)
# ------------------------------------------------------------------
$classif = $wb.Worksheets.Item("Classif Turno 2")

$ranking = @(
    @("Esquadrão Gazembrino", 83.5),
    @("SC 100 Sono", 67.16),
    @("GaúchoDaFronteira F.C", 63.76),
    @("bugredasmissões", 62.76),
    @("Texas Club 2026", 59.86),
    @("Doug Leal F.C", 58.17),
    @("Medonho´s F.C.", 55.96),
    @("lsauer fc", 54.36),
    @("Grêmio_Campeão_LA_27", 53.91),
    @("GrioTeam", 53.06),
    @("Pontaç0 F.C.", 52.66),
    @("Arran Katoko FC", 51.5),
    @("FBC Colorado", 45.46),
    @("GE Bebum", 44.65),
    @("La Primeira Patada Es Nuestra", 44.26),
    @("SC ÉoINTER!", 40.4),
    @("JV5 Tricolor Gaúcho", 39.66),
    @("C R Juvenal", 38.46),
    @("Pepe Leal FC", 16.4),
    @("NHU PORÃ SAF.", 0)
)

for ($i = 0; $i -lt $ranking.Length; $i++) {
    $row = $i + 2
    $classif.Cells.Item($row, 1).Value = $ranking[$i][0]
    $classif.Cells.Item($row, 2).Value = $ranking[$i][1]
}

Write-Output "edit applied"
